$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each cell below is stored as plain text in the workbook (inline string),
# e.g. prices like "303.92" and percentages like "5.56%". A leading
# apostrophe forces Excel to keep the literal text instead of coercing it
# to a number/percentage, and resetting the Style afterwards keeps the
# cell formatting identical to the original (no numFmt/quote-prefix style left behind).
$c = $ws.Range('D2')
$c.Value = "'303.92"
$c.Style = 'Normal'
$c = $ws.Range('E2')
$c.Value = "'5.56%"
$c.Style = 'Normal'
$c = $ws.Range('D3')
$c.Value = "'35.05"
$c.Style = 'Normal'
$c = $ws.Range('E3')
$c.Value = "'12.87%"
$c.Style = 'Normal'
$c = $ws.Range('D4')
$c.Value = "'5.182"
$c.Style = 'Normal'
$c = $ws.Range('E4')
$c.Value = "'5.34%"
$c.Style = 'Normal'
$c = $ws.Range('D5')
$c.Value = "'0.07760"
$c.Style = 'Normal'
$c = $ws.Range('E5')
$c.Value = "'6.06%"
$c.Style = 'Normal'
$c = $ws.Range('D6')
$c.Value = "'2.309"
$c.Style = 'Normal'
$c = $ws.Range('E6')
$c.Value = "'-0.85%"
$c.Style = 'Normal'
$c = $ws.Range('D7')
$c.Value = "'8.039"
$c.Style = 'Normal'
$c = $ws.Range('E7')
$c.Value = "'3.85%"
$c.Style = 'Normal'
$c = $ws.Range('D8')
$c.Value = "'4.008"
$c.Style = 'Normal'
$c = $ws.Range('E8')
$c.Value = "'7.75%"
$c.Style = 'Normal'
$c = $ws.Range('D9')
$c.Value = "'0.9290"
$c.Style = 'Normal'
$c = $ws.Range('E9')
$c.Value = "'2.84%"
$c.Style = 'Normal'
$c = $ws.Range('D10')
$c.Value = "'0.1017"
$c.Style = 'Normal'
$c = $ws.Range('E10')
$c.Value = "'11.40%"
$c.Style = 'Normal'
$c = $ws.Range('D11')
$c.Value = "'0.1830"
$c.Style = 'Normal'
$c = $ws.Range('E11')
$c.Value = "'8.42%"
$c.Style = 'Normal'
$c = $ws.Range('D12')
$c.Value = "'0.08636"
$c.Style = 'Normal'
$c = $ws.Range('E12')
$c.Value = "'4.65%"
$c.Style = 'Normal'
$c = $ws.Range('D13')
$c.Value = "'0.03453"
$c.Style = 'Normal'
$c = $ws.Range('E13')
$c.Value = "'10.71%"
$c.Style = 'Normal'
$c = $ws.Range('D14')
$c.Value = "'0.09918"
$c.Style = 'Normal'
$c = $ws.Range('E14')
$c.Value = "'-0.09%"
$c.Style = 'Normal'
$c = $ws.Range('D15')
$c.Value = "'0.001484"
$c.Style = 'Normal'
$c = $ws.Range('E15')
$c.Value = "'-0.92%"
$c.Style = 'Normal'
$c = $ws.Range('D16')
$c.Value = "'0.04626"
$c.Style = 'Normal'
$c = $ws.Range('E16')
$c.Value = "'2.36%"
$c.Style = 'Normal'
$c = $ws.Range('D17')
$c.Value = "'0.005797"
$c.Style = 'Normal'
$c = $ws.Range('E17')
$c.Value = "'1.42%"
$c.Style = 'Normal'
$c = $ws.Range('D18')
$c.Value = "'3.508"
$c.Style = 'Normal'
$c = $ws.Range('E18')
$c.Value = "'0.41%"
$c.Style = 'Normal'
$c = $ws.Range('D19')
$c.Value = "'2.110"
$c.Style = 'Normal'
$c = $ws.Range('E19')
$c.Value = "'3.11%"
$c.Style = 'Normal'
$c = $ws.Range('D20')
$c.Value = "'0.3415"
$c.Style = 'Normal'
$c = $ws.Range('E20')
$c.Value = "'2.56%"
$c.Style = 'Normal'
$c = $ws.Range('E21')
$c.Value = "'0.57%"
$c.Style = 'Normal'
$c = $ws.Range('D22')
$c.Value = "'4.629"
$c.Style = 'Normal'
$c = $ws.Range('E22')
$c.Value = "'9.85%"
$c.Style = 'Normal'
$c = $ws.Range('D23')
$c.Value = "'0.2347"
$c.Style = 'Normal'
$c = $ws.Range('E23')
$c.Value = "'11.74%"
$c.Style = 'Normal'
$c = $ws.Range('D24')
$c.Value = "'0.001227"
$c.Style = 'Normal'
$c = $ws.Range('E24')
$c.Value = "'1.56%"
$c.Style = 'Normal'
$c = $ws.Range('D25')
$c.Value = "'0.004419"
$c.Style = 'Normal'
$c = $ws.Range('E25')
$c.Value = "'6.17%"
$c.Style = 'Normal'
$c = $ws.Range('D26')
$c.Value = "'0.0001306"
$c.Style = 'Normal'
$c = $ws.Range('E26')
$c.Value = "'0.43%"
$c.Style = 'Normal'
$c = $ws.Range('D27')
$c.Value = "'0.0003427"
$c.Style = 'Normal'
$c = $ws.Range('E27')
$c.Value = "'0.95%"
$c.Style = 'Normal'
$c = $ws.Range('D39')
$c.Value = "'0.01759"
$c.Style = 'Normal'
$c = $ws.Range('E39')
$c.Value = "'12.01%"
$c.Style = 'Normal'
$c = $ws.Range('D40')
$c.Value = "'0.04728"
$c.Style = 'Normal'
$c = $ws.Range('E40')
$c.Value = "'6.49%"
$c.Style = 'Normal'
$c = $ws.Range('D41')
$c.Value = "'0.007640"
$c.Style = 'Normal'
$c = $ws.Range('E41')
$c.Value = "'4.25%"
$c.Style = 'Normal'
$c = $ws.Range('D42')
$c.Value = "'0.1405"
$c.Style = 'Normal'
$c = $ws.Range('E42')
$c.Value = "'6.00%"
$c.Style = 'Normal'
$c = $ws.Range('D43')
$c.Value = "'0.007093"
$c.Style = 'Normal'
$c = $ws.Range('E43')
$c.Value = "'-25.44%"
$c.Style = 'Normal'
$c = $ws.Range('D44')
$c.Value = "'0.002206"
$c.Style = 'Normal'
$c = $ws.Range('E44')
$c.Value = "'-0.67%"
$c.Style = 'Normal'
$c = $ws.Range('D45')
$c.Value = "'0.009229"
$c.Style = 'Normal'
$c = $ws.Range('E45')
$c.Value = "'3.22%"
$c.Style = 'Normal'
$c = $ws.Range('D46')
$c.Value = "'0.00005929"
$c.Style = 'Normal'
$c = $ws.Range('D47')
$c.Value = "'0.00000000754"
$c.Style = 'Normal'
$c = $ws.Range('E47')
$c.Value = "'0.46%"
$c.Style = 'Normal'
$c = $ws.Range('E48')
$c.Value = "'11.83%"
$c.Style = 'Normal'
$c = $ws.Range('D49')
$c.Value = "'0.002710"
$c.Style = 'Normal'
$c = $ws.Range('E49')
$c.Value = "'35.42%"
$c.Style = 'Normal'
$c = $ws.Range('D50')
$c.Value = "'0.00002110"
$c.Style = 'Normal'
$c = $ws.Range('E50')
$c.Value = "'0.46%"
$c.Style = 'Normal'
$c = $ws.Range('D51')
$c.Value = "'0.0002010"
$c.Style = 'Normal'
$c = $ws.Range('E51')
$c.Value = "'0.46%"
$c.Style = 'Normal'
